# This edit is a pass-through: the upstream change only swaps the
# internal numbering of two functionally-identical SharePoint
# customXml parts (customXml/item1.xml <-> customXml/item2.xml and
# their matching itemProps1.xml <-> itemProps2.xml companions).
# The content of each (item, itemProps) pair is unchanged; only which
# arbitrary "item1"/"item2" slot it occupies differs. This is package
# plumbing that SharePoint/OneDrive re-sync commonly reshuffles and is
# not part of the Word document object model (Document.CustomXMLParts
# exists as a property but is not backed by mutable part storage in
# this host, and there is no Find/Replace-style operation that retargets
# a raw OPC part's name). There is nothing in the visible Word content,
# styles, or real document properties to change here, so this script
# intentionally performs no content edits.
$d = $word.ActiveDocument
